$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (interested count) column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6444
$ws1.Range("F5").Value = 387
$ws1.Range("F9").Value = 88
$ws1.Range("F13").Value = 373
$ws1.Range("F14").Value = 943
$ws1.Range("F15").Value = 3152
$ws1.Range("F18").Value = 1822

# Sheet "全部类型" (all types) - same updates, rows shifted by +1 from row 8 onward
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6444
$ws4.Range("F5").Value = 387
$ws4.Range("F10").Value = 88
$ws4.Range("F14").Value = 373
$ws4.Range("F15").Value = 943
$ws4.Range("F16").Value = 3152
$ws4.Range("F19").Value = 1822
